$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 607.87
$ws.Range("I15").Value = 607.87
$ws.Range("K15").Value = 1823.61
$ws.Range("M15").Value = -1654.61
$ws.Range("H33").Value = 3637394.8
$ws.Range("I33").Value = 1199.9524
$ws.Range("J33").Value = 22727418
$ws.Range("K33").Value = 1199.9524
$ws.Range("L33").Value = 22727418
$ws.Range("M33").Value = -970.9523999999999
$ws.Range("N33").Value = -22727876
$ws.Range("H40").Value = 1776.0804
$ws.Range("I40").Value = 1811.1267
$ws.Range("J40").Value = 1620.5625
$ws.Range("K40").Value = 1811.1267
$ws.Range("L40").Value = 1620.5625
$ws.Range("M40").Value = -1636.1267
$ws.Range("N40").Value = -1970.5625
$ws.Range("H43").Value = 746
$ws.Range("I43").Value = 395
$ws.Range("J43").Value = 921.5
$ws.Range("K43").Value = 395
$ws.Range("L43").Value = 921.5
$ws.Range("M43").Value = -326
$ws.Range("N43").Value = -1059.5
$ws.Range("H137").Value = 2058.6897
$ws.Range("I137").Value = 1826.2632
$ws.Range("J137").Value = 2500.3
$ws.Range("K137").Value = 5478.7896
$ws.Range("L137").Value = 7500.900000000001
$ws.Range("M137").Value = -2928.7896
$ws.Range("N137").Value = -12600.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 635.34375
$ws.Range("I2").Value = 697.11536
$ws.Range("J2").Value = 367.66666
$ws.Range("K2").Value = 697.11536
$ws.Range("L2").Value = 367.66666
$ws.Range("M2").Value = -584.11536
$ws.Range("N2").Value = -593.66666
$ws.Range("H61").Value = 3339
$ws.Range("I61").Value = 3366.3953
$ws.Range("J61").Value = 2750
$ws.Range("K61").Value = 3366.3953
$ws.Range("L61").Value = 2750
$ws.Range("M61").Value = -3154.3953
$ws.Range("N61").Value = -3174
$ws.Range("H116").Value = 635.34375
$ws.Range("I116").Value = 697.11536
$ws.Range("J116").Value = 367.66666
$ws.Range("K116").Value = 697.11536
$ws.Range("L116").Value = 367.66666
$ws.Range("M116").Value = 1596.88464
$ws.Range("N116").Value = -4955.66666
$ws.Range("H132").Value = 1912.2073
$ws.Range("I132").Value = 1262.9403
$ws.Range("J132").Value = 4812.2666
$ws.Range("K132").Value = 3788.8209
$ws.Range("L132").Value = 14436.7998
$ws.Range("M132").Value = -1258.8209
$ws.Range("N132").Value = -19496.7998
$ws.Range("H136").Value = 3339
$ws.Range("I136").Value = 3366.3953
$ws.Range("J136").Value = 2750
$ws.Range("K136").Value = 10099.1859
$ws.Range("L136").Value = 8250
$ws.Range("M136").Value = -7549.1859
$ws.Range("N136").Value = -13350

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 635.34375
$ws.Range("I3").Value = 697.11536
$ws.Range("J3").Value = 367.66666
$ws.Range("K3").Value = 697.11536
$ws.Range("L3").Value = 367.66666
$ws.Range("M3").Value = -583.11536
$ws.Range("N3").Value = -595.66666
$ws.Range("H94").Value = 1406.5652
$ws.Range("I94").Value = 634.6923
$ws.Range("J94").Value = 2410
$ws.Range("K94").Value = 634.6923
$ws.Range("L94").Value = 2410
$ws.Range("M94").Value = -183.6923
$ws.Range("N94").Value = -3312
$ws.Range("H134").Value = 3036.2876
$ws.Range("I134").Value = 3057.3396
$ws.Range("J134").Value = 2980.5
$ws.Range("K134").Value = 9172.0188
$ws.Range("L134").Value = 8941.5
$ws.Range("M134").Value = -6637.0188
$ws.Range("N134").Value = -14011.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4418.3384
$ws.Range("I31").Value = 1595.119
$ws.Range("J31").Value = 8978.923000000001
$ws.Range("K31").Value = 1595.119
$ws.Range("L31").Value = 8978.923000000001
$ws.Range("M31").Value = -1300.119
$ws.Range("N31").Value = -9568.923000000001
$ws.Range("H34").Value = 4418.3384
$ws.Range("I34").Value = 1595.119
$ws.Range("J34").Value = 8978.923000000001
$ws.Range("K34").Value = 1595.119
$ws.Range("L34").Value = 8978.923000000001
$ws.Range("M34").Value = -1393.119
$ws.Range("N34").Value = -9382.923000000001
$ws.Range("H94").Value = 4509.96
$ws.Range("I94").Value = 4799.9
$ws.Range("J94").Value = 4316.6665
$ws.Range("K94").Value = 4799.9
$ws.Range("L94").Value = 4316.6665
$ws.Range("M94").Value = -4348.9
$ws.Range("N94").Value = -5218.6665
$ws.Range("H134").Value = 2066.2678
$ws.Range("I134").Value = 2154.0227
$ws.Range("J134").Value = 1744.5
$ws.Range("K134").Value = 6462.0681
$ws.Range("L134").Value = 5233.5
$ws.Range("M134").Value = -3927.0681
$ws.Range("N134").Value = -10303.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 62500348
$ws.Range("I26").Value = 166
$ws.Range("J26").Value = 100000456
$ws.Range("K26").Value = 498
$ws.Range("L26").Value = 300001368
$ws.Range("M26").Value = -210
$ws.Range("N26").Value = -300001944

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1075
$ws.Range("J97").Value = 1400
$ws.Range("L97").Value = 1400
$ws.Range("N97").Value = -2392
$ws.Range("H102").Value = 499575.9
$ws.Range("I102").Value = 606188.6
$ws.Range("J102").Value = 2050
$ws.Range("K102").Value = 606188.6
$ws.Range("L102").Value = 2050
$ws.Range("M102").Value = -604566.6
$ws.Range("N102").Value = -5294
$ws.Range("H113").Value = 38462744
$ws.Range("I113").Value = 76924020
$ws.Range("J113").Value = 1471.7693
$ws.Range("K113").Value = 76924020
$ws.Range("L113").Value = 1471.7693
$ws.Range("M113").Value = -76921850
$ws.Range("N113").Value = -5811.7693

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 680.6875
$ws.Range("I93").Value = 642.2143
$ws.Range("J93").Value = 950
$ws.Range("K93").Value = 642.2143
$ws.Range("L93").Value = 950
$ws.Range("M93").Value = 605.7857
$ws.Range("N93").Value = -3446
$ws.Range("H132").Value = 23238864
$ws.Range("I132").Value = 33404210
$ws.Range("J132").Value = 3788.5715
$ws.Range("K132").Value = 100212630
$ws.Range("L132").Value = 11365.7145
$ws.Range("M132").Value = -100210100
$ws.Range("N132").Value = -16425.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 20286.326
$ws.Range("I132").Value = 29442.486
$ws.Range("J132").Value = 1435.4117
$ws.Range("K132").Value = 88327.458
$ws.Range("L132").Value = 4306.2351
$ws.Range("M132").Value = -85797.458
$ws.Range("N132").Value = -9366.2351
